$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (new D value or $null, new E value or $null)
# D values that would otherwise be auto-parsed as a number are prefixed
# with a leading apostrophe so Excel stores them as literal text (matching
# the original inline-string cell content), exactly as typing '572.58 into
# a cell does interactively.
$updates = @{
    2  = @("63.066.57", "  +0.76%  ")
    3  = @("2.451.40", "  +0.87%  ")
    4  = @($null, "  -0.02%  ")
    5  = @("'572.58", "  +1.06%  ")
    6  = @("'146.15", "  +0.52%  ")
    7  = @($null, "  +0.07%  ")
    8  = @($null, "  +0.89%  ")
    9  = @("2.450.71", "  +0.89%  ")
    10 = @("'0.111", "  +0.65%  ")
    11 = @($null, "  +1.52%  ")
    12 = @("'5.25", "  -1.04%  ")
    13 = @($null, "  -0.01%  ")
    14 = @("'26.96", "  +0.61%  ")
    15 = @("'0.0000179", "  -0.58%  ")
    16 = @("2.901.45", "  +1.09%  ")
    17 = @("62.962.36", "  +0.77%  ")
    18 = @("2.447.78", $null)
    19 = @("'11.29", "  +0.53%  ")
    20 = @($null, "  +5.15%  ")
    21 = @("'328.31", "  +1.34%  ")
    22 = @($null, "  +0.89%  ")
    23 = @($null, "  +13.46%  ")
    24 = @($null, "  +0.75%  ")
    25 = @("'65.22", "  -2.92%  ")
    26 = @("'614.66", "  +2.60%  ")
    27 = @("'8.83", "  +2.94%  ")
    28 = @($null, "  +1.30%  ")
    29 = @("2.580.73", "  +1.21%  ")
    30 = @("'1.51", "  +4.04%  ")
    31 = @($null, "  +0.08%  ")
    32 = @($null, "  -2.79%  ")
    33 = @($null, "  +1.33%  ")
    34 = @($null, "  -1.87%  ")
    35 = @("'5.18", "  +6.51%  ")
    36 = @($null, "  +1.15%  ")
    37 = @($null, "  +0.08%  ")
    38 = @($null, "  -0.80%  ")
    39 = @("'18.85", "  +0.64%  ")
    40 = @($null, "  +0.42%  ")
    41 = @("'146.85", "  -0.18%  ")
    42 = @($null, "  -1.77%  ")
    43 = @("'2.59", "  +5.28%  ")
    44 = @($null, "  -0.12%  ")
    45 = @("'41.78", "  +0.42%  ")
    46 = @("'148.67", "  -0.01%  ")
    47 = @("'3.77", "  +2.46%  ")
    48 = @("'21.11", "  +2.62%  ")
    49 = @("'0.0532", "  -0.53%  ")
    50 = @("'0.602", "  +0.04%  ")
    51 = @($null, "  +0.68%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        $ws.Range("D$row").Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Range("E$row").Value = $eVal
    }
}
